$d = $word.ActiveDocument

$replacements = @(
    @("179×2=358", "896×2=1792"),
    @("275×7=1925", "215×4=860"),
    @("483×2=966", "159×4=636"),
    @("259×4=1036", "430×2=860"),
    @("867×8=6936", "737×6=4422"),
    @("952×2=1904", "216×2=432"),
    @("766×6=4596", "412×9=3708"),
    @("878×7=6146", "835×8=6680"),
    @("617×8=4936", "978×6=5868"),
    @("415×7=2905", "516×3=1548"),
    @("494×4=1976", "361×3=1083"),
    @("481×9=4329", "811×3=2433"),
    @("998×6=5988", "453×2=906"),
    @("939×7=6573", "161×5=805"),
    @("608×6=3648", "801×2=1602"),
    @("545×3=1635", "985×7=6895"),
    @("935×7=6545", "920×9=8280"),
    @("345×2=690", "363×3=1089"),
    @("382×4=1528", "414×7=2898"),
    @("185×9=1665", "224×6=1344"),
    @("620×6=3720", "493×4=1972"),
    @("160×5=800", "914×8=7312"),
    @("257×9=2313", "134×9=1206"),
    @("919×2=1838", "263×9=2367"),
    @("188×2=376", "698×7=4886"),
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}

$d.Save()
